# This script applies the "50% boolean" data update to the worksheet.
# It toggles the 0/1 boolean values in 97 cells across rows 3-11 (the
# GERMANY-LOW .. INDIA-HIGH data rows) to match the refreshed dataset,
# leaving row/column headers and all other cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 0
$ws.Range("V3").Value = 1
$ws.Range("Y3").Value = 0
$ws.Range("AC3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AJ3").Value = 1
$ws.Range("AL3").Value = 1
$ws.Range("AS3").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("T4").Value = 0
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 1
$ws.Range("Z4").Value = 1
$ws.Range("AA4").Value = 1
$ws.Range("AI4").Value = 1
$ws.Range("AR4").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 1
$ws.Range("W5").Value = 1
$ws.Range("X5").Value = 1
$ws.Range("AA5").Value = 1
$ws.Range("AC5").Value = 0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 1
$ws.Range("AI5").Value = 1
$ws.Range("AJ5").Value = 1
$ws.Range("AM5").Value = 0
$ws.Range("AN5").Value = 1
$ws.Range("AO5").Value = 1
$ws.Range("AS5").Value = 0
$ws.Range("AT5").Value = 1
$ws.Range("AU5").Value = 1
$ws.Range("AW5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("AB6").Value = 1
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 0
$ws.Range("AI6").Value = 1
$ws.Range("AK6").Value = 1
$ws.Range("AS6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("V7").Value = 1
$ws.Range("Y7").Value = 1
$ws.Range("Z7").Value = 0
$ws.Range("AL7").Value = 1
$ws.Range("AV7").Value = 1
$ws.Range("AX7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 1
$ws.Range("Q8").Value = 0
$ws.Range("W8").Value = 1
$ws.Range("Y8").Value = 1
$ws.Range("AI8").Value = 1
$ws.Range("AP8").Value = 0
$ws.Range("AX8").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("U9").Value = 1
$ws.Range("AA9").Value = 1
$ws.Range("AB9").Value = 0
$ws.Range("AC9").Value = 0
$ws.Range("AF9").Value = 0
$ws.Range("AL9").Value = 1
$ws.Range("AO9").Value = 0
$ws.Range("AP9").Value = 0
$ws.Range("AQ9").Value = 0
$ws.Range("AS9").Value = 0
$ws.Range("AT9").Value = 1
$ws.Range("AV9").Value = 1
$ws.Range("AW9").Value = 1
$ws.Range("AY9").Value = 0
$ws.Range("S10").Value = 1
$ws.Range("AA10").Value = 1
$ws.Range("AL10").Value = 1
$ws.Range("AP10").Value = 0
$ws.Range("AQ10").Value = 0
$ws.Range("AV10").Value = 1
$ws.Range("AQ11").Value = 0
